$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A112").Value = 'final_blood_smear_test'
$ws.Range("A113").Value = 'mothers_line_number'
$ws.Range("A114").Value = 'months_ago_net_obtained'
$ws.Range("A115").Value = 'malariae_present'
$ws.Range("A116").Value = 'blood_smear_bar_code'
$ws.Range("A117").Value = 'line_number'
$ws.Range("A118").Value = 'read_consent_statement_hemoglobin'
$ws.Range("A119").Value = 'ovale_present'
$ws.Range("A120").Value = 'someone_slept_under_net_last_night'
$ws.Range("A121").Value = 'day_of_birth'
$ws.Range("A122").Value = 'childs_age_in_months'
$ws.Range("A123").Value = 'net_design_no'
$ws.Range("A124").Value = 'na_read_consent_statement_for_malaria'
$ws.Range("A125").Value = 'falciparum_present'
$ws.Range("A126").Value = 'cmc_date_of_birth'
$ws.Range("A127").Value = 'month_of_data_collection'
$ws.Range("A128").Value = 'usual_resident'
$ws.Range("A129").Value = 'day_of_data_collection'
$ws.Range("A130").Value = 'malaria_measurement_result'
$ws.Range("A131").Value = 'children_hemoglobin_elig'
$ws.Range("A132").Value = 'hemoglobin_level_g_dl'
$ws.Range("A133").Value = 'childs_age_in_months_country_specific'
$ws.Range("A134").Value = 'childs_age_in_days_country_specific'
$ws.Range("A135").Value = 'century_day_code_of_birth'
$ws.Range("A136").Value = 'bed_net_type'
$ws.Range("A137").Value = 'child_age_in_months'
$ws.Range("A138").Value = 'index_to_household_schedule_hc0'
$ws.Range("A139").Value = 'childs_age_in_months_country_specific_hml16a'
$ws.Range("A140").Value = 'slept_llin_net'
$ws.Range("A141").Value = 'corr_age'
$ws.Range("A142").Value = 'result_of_measurement_hemoglobin'
$ws.Range("A143").Value = 'vivax_present'
$ws.Range("A144").Value = 'year_of_data_collection'
$ws.Range("A145").Value = 'date_measured_day'
$ws.Range("A146").Value = 'individual_file_pregnancy_status'
$ws.Range("A147").Value = 'month_of_birth'
$ws.Range("A148").Value = 'line_number_of_parent_caretaker'
$ws.Range("A149").Value = 'sex_of_member'
$ws.Range("A150").Value = 'slept_last_night'
$ws.Range("A151").Value = 'net_from_antenatal_immunization_visit'
$ws.Range("A152").Value = 'female_int_eligibility'
$ws.Range("A153").Value = 'number_of_persons_slept_under_net'
$ws.Range("A154").Value = 'sex'
$ws.Range("A155").Value = 'insecticide_treated_net'
$ws.Range("A156").Value = 'flag_age'
$ws.Range("A157").Value = 'rapid_test_result'
$ws.Range("A158").Value = 'date_measured_month'
$ws.Range("A159").Value = 'slept_under_net'
$ws.Range("A160").Value = 'childs_age_in_days'
$ws.Range("A161").Value = 'completeleness_of_hc32_info'
$ws.Range("A162").Value = 'caretaker_line_number'
$ws.Range("A163").Value = 'index_to_household_schedule_hmhidx'
$ws.Range("A164").Value = 'fieldworker_measurer_code'
$ws.Range("A165").Value = 'hemoglobin_level_adjusted_for_altitude_g_dl'
$ws.Range("A166").Value = 'year_of birth'
$ws.Range("A167").Value = 'anemia_level'
$ws.Range("A168").Value = 'century_day_code_of_measurement'
$ws.Range("A169").Value = 'rshp_to_head'
$ws.Range("A170").Value = 'date_measured_year'
$ws.Range("A171").Value = 'mosquito_bed_net_designation_number'
$ws.Range("A172").Value = 'line_number_of_person_slept_in_net'
$ws.Range("A173").Value = 'net_observed_by_interviewer'
$ws.Range("A174").Value = 'brand_of_net'
$ws.Range("A175").Value = 'age_of_member'
$ws.Range("A176").Value = 'fieldworker_malaria_measurer_code'
